# Update "想去人数" (number of people wanting to go) values in column F
# across the four worksheets, per the source diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1656
$ws1.Range("F3").Value = 864
$ws1.Range("F4").Value = 274
$ws1.Range("F6").Value = 1188
$ws1.Range("F7").Value = 807
$ws1.Range("F8").Value = 830
$ws1.Range("F9").Value = 1533
$ws1.Range("F11").Value = 1062
$ws1.Range("F12").Value = 33
$ws1.Range("F13").Value = 76
$ws1.Range("F14").Value = 204
$ws1.Range("F16").Value = 512
$ws1.Range("F19").Value = 10
$ws1.Range("F24").Value = 55
$ws1.Range("F25").Value = 9
$ws1.Range("F27").Value = 262

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1037
$ws2.Range("F5").Value = 283
$ws2.Range("F8").Value = 72
$ws2.Range("F9").Value = 599
$ws2.Range("F10").Value = 91

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 268

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 268
$ws4.Range("F3").Value = 1656
$ws4.Range("F5").Value = 864
$ws4.Range("F6").Value = 274
$ws4.Range("F7").Value = 1037
$ws4.Range("F9").Value = 1188
$ws4.Range("F10").Value = 807
$ws4.Range("F11").Value = 830
$ws4.Range("F12").Value = 1533
$ws4.Range("F14").Value = 1062
$ws4.Range("F15").Value = 33
$ws4.Range("F16").Value = 76
$ws4.Range("F17").Value = 204
$ws4.Range("F19").Value = 512
$ws4.Range("F23").Value = 10
$ws4.Range("F24").Value = 283
$ws4.Range("F32").Value = 55
$ws4.Range("F33").Value = 9
$ws4.Range("F35").Value = 262
$ws4.Range("F36").Value = 72
$ws4.Range("F38").Value = 599
$ws4.Range("F39").Value = 91
$ws4.Range("F40").Value = 91
